$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "60.249.81"
    "E2" = "  -1.26%  "
    "D3" = "2.588.67"
    "E3" = "  -3.33%  "
    "E4" = "  +0.39%  "
    "D5" = "508.56"
    "E5" = "  -0.97%  "
    "D6" = "153.38"
    "E6" = "  -2.73%  "
    "D7" = "0.998"
    "E7" = "  +0.30%  "
    "D8" = "0.588"
    "E8" = "  -4.08%  "
    "D9" = "2.596.05"
    "E9" = "  -3.00%  "
    "E10" = "  +6.37%  "
    "D11" = "0.103"
    "E11" = "  -1.85%  "
    "D12" = "0.345"
    "E12" = "  -1.00%  "
    "E13" = "  +1.66%  "
    "D14" = "3.045.08"
    "E14" = "  -1.85%  "
    "D15" = "60.236.57"
    "E15" = "  -1.32%  "
    "D16" = "21.50"
    "E16" = "  -1.85%  "
    "E17" = "  -0.45%  "
    "D18" = "2.596.72"
    "E18" = "  -2.82%  "
    "E19" = "  -1.36%  "
    "D20" = "352.43"
    "E20" = "  +0.36%  "
    "D21" = "10.50"
    "E21" = "  -0.23%  "
    "E22" = "  -1.29%  "
    "E23" = "  +0.04%  "
    "D24" = "60.32"
    "E24" = "  +0.16%  "
    "E25" = "  -0.76%  "
    "E26" = "  -0.36%  "
    "E27" = "  +0.90%  "
    "D28" = "0.0₃0835"
    "E28" = "  -4.09%  "
    "D29" = "7.32"
    "E29" = "  -3.12%  "
    "E30" = "  +0.34%  "
    "D31" = "19.36"
    "D32" = "151.58"
    "E32" = "  -3.79%  "
    "E33" = "  -1.53%  "
    "E34" = "  -0.08%  "
    "D35" = "3.98"
    "E35" = "  -1.59%  "
    "E36" = "  -3.26%  "
    "E37" = "  +3.53%  "
    "D38" = "1.47"
    "E38" = "  -3.01%  "
    "D39" = "36.08"
    "E39" = "  +1.84%  "
    "B40" = "Filecoin"
    "C40" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D40" = "3.75"
    "E40" = "  -0.89%  "
    "B41" = "Fetch.AI"
    "C41" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D41" = "0.837"
    "E41" = "  -3.72%  "
    "D42" = "295.01"
    "E42" = "  -4.95%  "
    "D43" = "0.100"
    "E43" = "  -1.19%  "
    "D44" = "0.616"
    "E44" = "  -4.81%  "
    "E45" = "  -0.20%  "
    "D46" = "0.0550"
    "E46" = "  -4.64%  "
    "D47" = "19.65"
    "E47" = "  -2.03%  "
    "E48" = "  -5.03%  "
    "D49" = "0.0232"
    "E49" = "  -1.97%  "
    "E50" = "  -0.18%  "
    "D51" = "1.987.86"
    "E51" = "  -2.55%  "
}

foreach ($cellRef in $updates.Keys) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "21.50", "0.0550")
    # keep their exact formatting instead of being auto-coerced to a Number.
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cellRef]
}

foreach ($cellRef in $updates.Keys) {
    # Drop the temporary text-format override so the cell style matches the original (no explicit "s").
    $ws.Range($cellRef).Style = "Normal"
}
